# Refresh the cryptocurrency price / 1h-volume-change table with the
# latest scraped values (GitHub Actions data refresh).
#
# Rows 19 and 20 additionally swap which coin (Uniswap / WrappedBTC)
# occupies each rank position, per the upstream source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.240.18"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").Value = "3.385.78"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.01"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.02"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +2.09%  "

$ws.Range("D8").Value = "3.373.86"
$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +5.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.638"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.76"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.19"
$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("D15").Value = "3.921.38"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.33"
$ws.Range("E16").Value = "  -0.64%  "

$ws.Range("D17").Value = "3.395.53"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "65.201.70"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.88"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.18"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.91"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.40"
$ws.Range("E24").Value = "  +8.37%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.39"
$ws.Range("E26").Value = "  +2.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.73"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.16"
$ws.Range("E30").Value = "  +3.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.55"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "63.08"
$ws.Range("E32").Value = "  +7.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.48"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "580.95"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  +4.83%  "

$ws.Range("E38").Value = "  +1.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.77"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("E41").Value = "  -2.16%  "

$ws.Range("D42").Value = "3.092.89"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0417"
$ws.Range("E43").Value = "  +1.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("E46").Value = "  -3.33%  "

$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.75"
$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.32"
$ws.Range("E51").Value = "  -0.61%  "
